# Update 2021 HWL2 First Batch
# Adds six new indicator rows (79-84) to Sheet1, gives the first new row's
# label cell a wrapped-text style, nudges a couple of cosmetic view
# settings (zoom / selection / tab ratio) to match the refreshed workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New rows of data (IndicatorName, LogScale, GraphType)
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 79; Name = "Global Extreme Poverty Cost of Basic Needs"; LogScale = 0; GraphType = "boxplot" },
    @{ Row = 80; Name = "Global Extreme Poverty Dollar a Day ";       LogScale = 0; GraphType = "boxplot" },
    @{ Row = 81; Name = "Wealth Yearly Ginis";                        LogScale = 0; GraphType = "boxplot" },
    @{ Row = 82; Name = "Wealth Total";                               LogScale = 0; GraphType = "boxplot" },
    @{ Row = 83; Name = "Wealth Top10 percent share";                 LogScale = 0; GraphType = "boxplot" },
    @{ Row = 84; Name = "Wealth Decadal Ginis";                       LogScale = 0; GraphType = "boxplot" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Name
    $ws.Cells.Item($rowNum, 2).Value = $r.LogScale
    $ws.Cells.Item($rowNum, 3).Value = $r.GraphType
}

# First new indicator name wraps onto two lines in the source column width.
$ws.Range("A79").WrapText = $true

# Row heights picked up by the refreshed sheet (wrapped / taller rows).
$ws.Rows.Item(79).RowHeight = 12.85
$ws.Rows.Item(80).RowHeight = 14.65
$ws.Rows.Item(81).RowHeight = 12.8
$ws.Rows.Item(82).RowHeight = 12.8
$ws.Rows.Item(83).RowHeight = 12.8
$ws.Rows.Item(84).RowHeight = 12.8

# ---------------------------------------------------------------------
# View / window state
# ---------------------------------------------------------------------
$ws.Range("B84:C84").Select()
$excel.ActiveWindow.Zoom = 180
$excel.ActiveWindow.ScrollRow = 68
$excel.ActiveWindow.ScrollColumn = 1

# Tab ratio of the split between sheet tabs and the horizontal scrollbar.
$excel.ActiveWindow.TabRatio = 500
